$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 6670776
$ws.Cells.Item(86, 9).Value = 14289549
$ws.Cells.Item(86, 10).Value = 4349.875
$ws.Cells.Item(86, 11).Value = 14289549
$ws.Cells.Item(86, 12).Value = 4349.875
$ws.Cells.Item(86, 13).Value = -14288426
$ws.Cells.Item(86, 14).Value = -6595.875
$ws.Cells.Item(89, 8).Value = 6670776
$ws.Cells.Item(89, 9).Value = 14289549
$ws.Cells.Item(89, 10).Value = 4349.875
$ws.Cells.Item(89, 11).Value = 71447745
$ws.Cells.Item(89, 12).Value = 21749.375
$ws.Cells.Item(89, 13).Value = -71442129
$ws.Cells.Item(89, 14).Value = -32981.375
$ws.Cells.Item(113, 8).Value = 2164.1428
$ws.Cells.Item(113, 9).Value = 1839.6
$ws.Cells.Item(113, 10).Value = 2344.4443
$ws.Cells.Item(113, 11).Value = 1839.6
$ws.Cells.Item(113, 12).Value = 2344.4443
$ws.Cells.Item(113, 13).Value = 1414.4
$ws.Cells.Item(113, 14).Value = -8852.444299999999
$ws.Cells.Item(116, 8).Value = 3571.2917
$ws.Cells.Item(116, 9).Value = 2987.3333
$ws.Cells.Item(116, 10).Value = 4544.5557
$ws.Cells.Item(116, 11).Value = 2987.3333
$ws.Cells.Item(116, 12).Value = 4544.5557
$ws.Cells.Item(116, 13).Value = 454.6667000000002
$ws.Cells.Item(116, 14).Value = -11428.5557
$ws.Cells.Item(129, 8).Value = 1152.0227
$ws.Cells.Item(129, 9).Value = 1515.6666
$ws.Cells.Item(129, 11).Value = 4546.9998
$ws.Cells.Item(129, 13).Value = 453.0002000000004
$ws.Cells.Item(138, 8).Value = 1372.02
$ws.Cells.Item(138, 9).Value = 620.88635
$ws.Cells.Item(138, 10).Value = 1962.1964
$ws.Cells.Item(138, 11).Value = 1862.65905
$ws.Cells.Item(138, 12).Value = 5886.5892
$ws.Cells.Item(138, 13).Value = 3277.34095
$ws.Cells.Item(138, 14).Value = -16166.5892

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2094.75
$ws.Cells.Item(45, 9).Value = 1930.2778
$ws.Cells.Item(45, 11).Value = 1930.2778
$ws.Cells.Item(45, 13).Value = -1553.2778
$ws.Cells.Item(61, 8).Value = 1125.6586
$ws.Cells.Item(61, 9).Value = 850.9666999999999
$ws.Cells.Item(61, 11).Value = 850.9666999999999
$ws.Cells.Item(61, 13).Value = -638.9666999999999
$ws.Cells.Item(88, 8).Value = 24231408
$ws.Cells.Item(88, 9).Value = 200000000
$ws.Cells.Item(88, 10).Value = 4701563.5
$ws.Cells.Item(88, 11).Value = 200000000
$ws.Cells.Item(88, 12).Value = 4701563.5
$ws.Cells.Item(88, 13).Value = -199999594
$ws.Cells.Item(88, 14).Value = -4702375.5
$ws.Cells.Item(91, 8).Value = 24231408
$ws.Cells.Item(91, 9).Value = 200000000
$ws.Cells.Item(91, 10).Value = 4701563.5
$ws.Cells.Item(91, 11).Value = 200000000
$ws.Cells.Item(91, 12).Value = 4701563.5
$ws.Cells.Item(91, 13).Value = -199998596
$ws.Cells.Item(91, 14).Value = -4704371.5
$ws.Cells.Item(115, 8).Value = 19999.889
$ws.Cells.Item(115, 10).Value = 19999.889
$ws.Cells.Item(115, 12).Value = 19999.889
$ws.Cells.Item(115, 14).Value = -23133.889
$ws.Cells.Item(123, 8).Value = 35000
$ws.Cells.Item(123, 10).Value = 35000
$ws.Cells.Item(123, 12).Value = 35000
$ws.Cells.Item(123, 14).Value = -44800
$ws.Cells.Item(136, 8).Value = 1125.6586
$ws.Cells.Item(136, 9).Value = 850.9666999999999
$ws.Cells.Item(136, 11).Value = 2552.9001
$ws.Cells.Item(136, 13).Value = -2.900099999999838

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1251849.9
$ws.Cells.Item(86, 9).Value = 2200.5
$ws.Cells.Item(86, 10).Value = 2501499.2
$ws.Cells.Item(86, 11).Value = 2200.5
$ws.Cells.Item(86, 12).Value = 2501499.2
$ws.Cells.Item(86, 13).Value = -1077.5
$ws.Cells.Item(86, 14).Value = -2503745.2
$ws.Cells.Item(89, 8).Value = 1251849.9
$ws.Cells.Item(89, 9).Value = 2200.5
$ws.Cells.Item(89, 10).Value = 2501499.2
$ws.Cells.Item(89, 11).Value = 11002.5
$ws.Cells.Item(89, 12).Value = 12507496
$ws.Cells.Item(89, 13).Value = -5386.5
$ws.Cells.Item(89, 14).Value = -12518728
$ws.Cells.Item(94, 8).Value = 774.03925
$ws.Cells.Item(94, 9).Value = 816.9773
$ws.Cells.Item(94, 10).Value = 504.14285
$ws.Cells.Item(94, 11).Value = 816.9773
$ws.Cells.Item(94, 12).Value = 504.14285
$ws.Cells.Item(94, 13).Value = -365.9773
$ws.Cells.Item(94, 14).Value = -1406.14285
$ws.Cells.Item(124, 8).Value = 43593.332
$ws.Cells.Item(124, 10).Value = 43593.332
$ws.Cells.Item(124, 12).Value = 43593.332
$ws.Cells.Item(124, 14).Value = -53413.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4392.602
$ws.Cells.Item(31, 9).Value = 1953.6666
$ws.Cells.Item(31, 10).Value = 5019.7573
$ws.Cells.Item(31, 11).Value = 1953.6666
$ws.Cells.Item(31, 12).Value = 5019.7573
$ws.Cells.Item(31, 13).Value = -1658.6666
$ws.Cells.Item(31, 14).Value = -5609.7573
$ws.Cells.Item(34, 8).Value = 4392.602
$ws.Cells.Item(34, 9).Value = 1953.6666
$ws.Cells.Item(34, 10).Value = 5019.7573
$ws.Cells.Item(34, 11).Value = 1953.6666
$ws.Cells.Item(34, 12).Value = 5019.7573
$ws.Cells.Item(34, 13).Value = -1751.6666
$ws.Cells.Item(34, 14).Value = -5423.7573
$ws.Cells.Item(132, 8).Value = 120867.25
$ws.Cells.Item(132, 9).Value = 1259.6
$ws.Cells.Item(132, 10).Value = 206301.28
$ws.Cells.Item(132, 11).Value = 3778.8
$ws.Cells.Item(132, 12).Value = 618903.84
$ws.Cells.Item(132, 13).Value = -1248.8
$ws.Cells.Item(132, 14).Value = -623963.84
$ws.Cells.Item(134, 8).Value = 453218.03
$ws.Cells.Item(134, 9).Value = 1159.75
$ws.Cells.Item(134, 10).Value = 2003132.1
$ws.Cells.Item(134, 11).Value = 3479.25
$ws.Cells.Item(134, 12).Value = 6009396.300000001
$ws.Cells.Item(134, 13).Value = -944.25
$ws.Cells.Item(134, 14).Value = -6014466.300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 5114.7144
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 5114.7144
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 15344.1432
$ws.Cells.Item(68, 14).Value = -16966.1432
$ws.Cells.Item(68, 13).ClearContents()
$ws.Cells.Item(71, 8).Value = 5114.7144
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 5114.7144
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 46032.4296
$ws.Cells.Item(71, 14).Value = -54144.4296
$ws.Cells.Item(71, 13).ClearContents()
$ws.Cells.Item(107, 8).Value = 6837.7095
$ws.Cells.Item(107, 9).Value = 9543
$ws.Cells.Item(107, 10).Value = 5349.8
$ws.Cells.Item(107, 11).Value = 28629
$ws.Cells.Item(107, 12).Value = 16049.4
$ws.Cells.Item(107, 13).Value = -26709
$ws.Cells.Item(107, 14).Value = -19889.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4624.857
$ws.Cells.Item(70, 9).Value = 4645.6665
$ws.Cells.Item(70, 10).Value = 4500
$ws.Cells.Item(70, 11).Value = 4645.6665
$ws.Cells.Item(70, 12).Value = 4500
$ws.Cells.Item(70, 13).Value = -4375.6665
$ws.Cells.Item(70, 14).Value = -5040
$ws.Cells.Item(73, 8).Value = 4624.857
$ws.Cells.Item(73, 9).Value = 4645.6665
$ws.Cells.Item(73, 10).Value = 4500
$ws.Cells.Item(73, 11).Value = 4645.6665
$ws.Cells.Item(73, 12).Value = 4500
$ws.Cells.Item(73, 13).Value = -3709.6665
$ws.Cells.Item(73, 14).Value = -6372
$ws.Cells.Item(75, 8).Value = 39500
$ws.Cells.Item(75, 10).Value = 39500
$ws.Cells.Item(75, 12).Value = 39500
$ws.Cells.Item(75, 14).Value = -41248
$ws.Cells.Item(78, 8).Value = 39500
$ws.Cells.Item(78, 10).Value = 39500
$ws.Cells.Item(78, 12).Value = 118500
$ws.Cells.Item(78, 14).Value = -127236
$ws.Cells.Item(97, 8).Value = 1926.4348
$ws.Cells.Item(97, 9).Value = 2014.6666
$ws.Cells.Item(97, 10).Value = 1000
$ws.Cells.Item(97, 11).Value = 2014.6666
$ws.Cells.Item(97, 12).Value = 1000
$ws.Cells.Item(97, 13).Value = -1518.6666
$ws.Cells.Item(97, 14).Value = -1992
$ws.Cells.Item(132, 8).Value = 3169.4138
$ws.Cells.Item(132, 9).Value = 2219.0588
$ws.Cells.Item(132, 10).Value = 4515.75
$ws.Cells.Item(132, 11).Value = 6657.176399999999
$ws.Cells.Item(132, 12).Value = 13547.25
$ws.Cells.Item(132, 13).Value = -4127.176399999999
$ws.Cells.Item(132, 14).Value = -18607.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 1771
$ws.Cells.Item(136, 9).Value = 1559.65
$ws.Cells.Item(136, 10).Value = 2155.2727
$ws.Cells.Item(136, 11).Value = 4678.950000000001
$ws.Cells.Item(136, 12).Value = 6465.8181
$ws.Cells.Item(136, 13).Value = -2128.950000000001
$ws.Cells.Item(136, 14).Value = -11565.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 703
$ws.Cells.Item(96, 9).Value = 703
$ws.Cells.Item(96, 11).Value = 703
$ws.Cells.Item(96, 13).Value = 670
$ws.Cells.Item(107, 8).Value = 901.44446
$ws.Cells.Item(107, 9).Value = 775
$ws.Cells.Item(107, 10).Value = 1002.6
$ws.Cells.Item(107, 11).Value = 2325
$ws.Cells.Item(107, 12).Value = 3007.8
$ws.Cells.Item(107, 13).Value = -405
$ws.Cells.Item(107, 14).Value = -6847.8
$ws.Cells.Item(123, 8).Value = 49106
$ws.Cells.Item(123, 10).Value = 49106
$ws.Cells.Item(123, 12).Value = 49106
$ws.Cells.Item(123, 14).Value = -58906
